# Refresh the crypto price/volume snapshot (cryptos list), mirroring the
# data pulled by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") and column E ("Volume(1h)") are stored as plain text in
# this sheet (prices use "." as a thousands separator, e.g. "69.761.52", and
# percentages keep padding spaces), not as numbers. Whenever a new price looks
# like a genuine number (e.g. "574.45"), force the cell to Text format first so
# Excel keeps storing it as a string instead of silently converting it.

$ws.Range("D2").Value = "70.026.19"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "3.559.12"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.45"
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.31"
$ws.Range("E6").Value = "  -2.75%  "

$ws.Range("D7").Value = "3.555.52"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.618"
$ws.Range("E8").Value = "  -2.15%  "

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.187"
$ws.Range("E10").Value = "  +5.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.646"
$ws.Range("E11").Value = "  -2.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.32"
$ws.Range("E12").Value = "  -4.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000302"
$ws.Range("E13").Value = "  +0.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.57"
$ws.Range("E14").Value = "  -2.66%  "

$ws.Range("D15").Value = "4.124.27"
$ws.Range("E15").Value = "  -0.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.49"
$ws.Range("E16").Value = "  -3.27%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "69.943.46"
$ws.Range("E17").Value = "  +0.33%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.535.25"
$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.47"
$ws.Range("E19").Value = "  -0.50%  "

$ws.Range("E20").Value = "  -1.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.05"
$ws.Range("E21").Value = "  +0.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "532.03"
$ws.Range("E22").Value = "  +11.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.68"
$ws.Range("E23").Value = "  +0.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.94"
$ws.Range("E24").Value = "  -2.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.40"
$ws.Range("E25").Value = "  +0.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.59"
$ws.Range("E26").Value = "  +7.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.29"
$ws.Range("E27").Value = "  +1.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.95"
$ws.Range("E28").Value = "  -4.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.24"
$ws.Range("E29").Value = "  -0.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.74"
$ws.Range("E30").Value = "  -1.28%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.53"
$ws.Range("E31").Value = "  -2.78%  "

$ws.Range("E32").Value = "  +2.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "65.49"
$ws.Range("E33").Value = "  -1.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.116"
$ws.Range("E34").Value = "  -4.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "578.09"
$ws.Range("E35").Value = "  -3.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.19"
$ws.Range("E36").Value = "  +6.92%  "

$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.406"
$ws.Range("E37").Value = "  +1.31%  "

$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.46"
$ws.Range("E38").Value = "  -5.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("D40").Value = "0.0₃0779"
$ws.Range("E40").Value = "  -3.58%  "

$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.18"
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.43"
$ws.Range("E42").Value = "  -2.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.135"
$ws.Range("E43").Value = "  -7.07%  "

$ws.Range("E44").Value = "  +6.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.01"
$ws.Range("E45").Value = "  -3.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0449"
$ws.Range("E46").Value = "  +0.56%  "

$ws.Range("D47").Value = "3.180.05"
$ws.Range("E47").Value = "  -2.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.42"
$ws.Range("E48").Value = "  -1.28%  "

$ws.Range("E49").Value = "  -1.44%  "

$ws.Range("E50").Value = "  +25.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.995"
$ws.Range("E51").Value = "  -0.40%  "
